# "added low grades back in with error bars"
# Append a new work-log entry (row 34) to Sheet1: 2024-06-12 (Wed), 4 hours,
# note about re-including low grades / quantile error bars.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 34

$ws.Cells.Item($newRow, 1).Value = 45455
$ws.Cells.Item($newRow, 2).Value = "W"
$ws.Cells.Item($newRow, 3).Value = 4
$ws.Cells.Item($newRow, 5).Value = "Still running grid search, added back in grades < 60. Removing rows without prereqs helped a bunch on train/test. Quantile error bars look good w new data"

# Match the date/day-of-week formatting and wrapped-note formatting used by
# the rest of the log.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item($newRow - 1, 2).NumberFormat
$ws.Cells.Item($newRow, 5).WrapText = $true

# Same wrapped-note row height as other multi-line entries.
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item(27).RowHeight

$ws.Range("C34").Select()
